# "Added last meeting minutes"
#
# Near the end of the document, after the "Any other issues or comments?"
# heading paragraph, there is one blank paragraph followed by six more
# blank paragraphs that are left over formatting cruft, followed by a run
# of paragraphs that carry the minutes-style indent (ind left="697"
# hanging="697"). The edit:
#   1. removes the six superfluous blank paragraphs, and
#   2. resets the indent of the paragraph that used to follow them
#      (now the first paragraph of the minutes block) to
#      ind left="0" firstLine="0", so it lines up flush-left as the
#      start of the newly added minutes content.

$d = $word.ActiveDocument

# Find the anchor heading paragraph.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Any other issues or comments?*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not locate anchor paragraph 'Any other issues or comments?'"
}

# One blank paragraph right after the heading is kept as-is.
$keepIndex = $anchorIndex + 1

# The next six blank paragraphs are removed entirely.
$deleteStart = $keepIndex + 1
$deleteEnd = $deleteStart + 5

$startPara = $d.Paragraphs.Item($deleteStart)
$endPara = $d.Paragraphs.Item($deleteEnd)
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()

# The paragraph that used to sit right after the deleted block (and now
# immediately follows the kept blank paragraph) gets its hanging indent
# replaced with a flush-left, no-first-line indent.
$target = $d.Paragraphs.Item($deleteStart)
$target.Range.ParagraphFormat.LeftIndent = 0
$target.Range.ParagraphFormat.FirstLineIndent = 0
